# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, matching the new scrape output.

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F (both sheets share identical data).
$updates = @{
    3  = 3230
    4  = 237
    5  = 132
    6  = 203
    7  = 1700
    8  = 1638
    9  = 470
    10 = 375
    13 = 192
    15 = 233
    21 = 57
    23 = 381
    24 = 226
    26 = 38
    27 = 11
    28 = 27
    29 = 309
    30 = 2208
    33 = 473
    34 = 338
    36 = 427
    37 = 227
    39 = 414
    40 = 521
    41 = 415
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
